$wb = $excel.ActiveWorkbook

# --- Rename existing "Requested quantity" headers ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$data = @(
    @(44941.99999999999, 136, -1019.710561551587, 1208.142410817111),
    @(44948.99999999999, 160, -923.8386116780616, 1264.281285489657),
    @(44955.99999999999, 183, -811.1700320671538, 1285.569360024796),
    @(44962.99999999999, 206, -898.0877855225391, 1343.899557685766),
    @(44976.99999999999, 253, -781.5075432704537, 1363.78921807806),
    @(44990.99999999999, 300, -796.1526197583029, 1369.952470827086),
    @(44997.99999999999, 324, -733.6881176857344, 1428.836863139327),
    @(45004.99999999999, 347, -686.5377827279358, 1411.150875724294),
    @(45018.99999999999, 394, -706.2651036296493, 1450.774707883298),
    @(45025.99999999999, 417, -561.369892059844, 1436.886759059493),
    @(45032.99999999999, 441, -546.4130556011286, 1512.416307196695),
    @(45039.99999999999, 464, -608.4590547113453, 1443.950710123981),
    @(45046.99999999999, 488, -605.0666749460955, 1547.834382943766),
    @(45053.99999999999, 511, -589.0269537892642, 1622.679417707264),
    @(45060.99999999999, 534, -581.6237158962268, 1566.0885273125),
    @(45067.99999999999, 558, -496.3433892326765, 1613.090997183738),
    @(45088.99999999999, 628, -433.0798608181824, 1748.640175554808),
    @(45095.99999999999, 652, -423.1402418773943, 1696.437079844617),
    @(45102.99999999999, 675, -470.7581577751814, 1679.026351590082),
    @(45109.99999999999, 699, -269.8390263352674, 1778.566805547465),
    @(45116.99999999999, 722, -274.5501566806174, 1844.441678678397),
    @(45123.99999999999, 745, -365.1587773600219, 1825.668673990429),
    @(45130.99999999999, 769, -301.9464805713115, 1917.965301257516),
    @(45144.99999999999, 816, -179.5908542251864, 1913.814646261444),
    @(45158.99999999999, 863, -212.083035157442, 1934.678470892075),
    @(45165.99999999999, 886, -210.6968322237452, 1856.577066504712),
    @(45179.99999999999, 933, -133.9163527571277, 2039.227340888072),
    @(45186.99999999999, 956, -111.6570486859609, 2055.334950027684),
    @(45200.99999999999, 1003, -82.81016162297871, 2119.999379752697),
    @(45207.99999999999, 1027, -29.2558961968464, 2015.12481827106),
    @(45214.99999999999, 1050, 51.16267305871173, 2103.362694790033),
    @(45221.99999999999, 1073, -41.44524747931837, 2126.856069340157),
    @(45249.99999999999, 1167, 60.03396507315716, 2220.390173142905),
    @(45256.99999999999, 1191, 111.0870424414452, 2327.885966020684),
    @(45263.99999999999, 1214, 159.5635766073236, 2328.387228583072),
    @(45270.99999999999, 1237, 245.5999198851346, 2323.525357089067),
    @(45277.99999999999, 1261, 210.5945100645924, 2297.351434876816),
    @(45298.99999999999, 1331, 308.783446092827, 2391.987713470199),
    @(45305.99999999999, 1355, 339.0572928542383, 2463.21857287062),
    @(45312.99999999999, 1378, 270.7723414993004, 2484.092086283807),
    @(45319.99999999999, 1402, 358.3940618588675, 2354.698329143563),
    @(45326.99999999999, 1425, 266.9551888471149, 2456.78912073663),
    @(45333.99999999999, 1448, 407.55506823853, 2508.107696876227),
    @(45340.99999999999, 1472, 363.4502129535528, 2530.047671091711),
    @(45347.99999999999, 1495, 346.6856471598484, 2621.513763330339),
    @(45354.99999999999, 1519, 410.2489276783133, 2636.060132768826),
    @(45361.99999999999, 1542, 495.4784660207414, 2648.092111875223),
    @(45368.99999999999, 1566, 464.078349548729, 2638.005590842937),
    @(45375.99999999999, 1589, 548.3906254607745, 2656.416074700501),
    @(45382.99999999999, 1612, 518.948410794071, 2766.249945660742),
    @(45389.99999999999, 1636, 490.157285330218, 2753.791713870463),
    @(45396.99999999999, 1659, 570.2493164914935, 2646.920742481278),
    @(45403.99999999999, 1683, 576.1185222495171, 2760.027326587056),
    @(45410.99999999999, 1706, 623.2336678990733, 2766.390083962171),
    @(45417.99999999999, 1730, 625.479166175606, 2741.053569804774),
    @(45424.99999999999, 1753, 653.0999613540989, 2816.845852693689),
    @(45431.99999999999, 1776, 731.3067825200661, 2854.357120639715),
    @(45438.99999999999, 1800, 624.1717999501366, 2855.333955384375),
    @(45445.99999999999, 1823, 784.985300606537, 2899.230162258582),
    @(45452.99999999999, 1847, 728.0335711652308, 2908.164653392329),
    @(45459.99999999999, 1870, 854.4958023558954, 2936.777874017107),
    @(45466.99999999999, 1894, 855.4319371148114, 2992.695146353601),
    @(45473.99999999999, 1917, 813.2829867439596, 3048.356713955393),
    @(45487.99999999999, 1964, 913.4786929544647, 3039.524731329479),
    @(45494.99999999999, 1987, 900.409576594041, 3044.547196753142),
    @(45501.99999999999, 2011, 971.4174366419354, 3038.079333969632),
    @(45508.99999999999, 2034, 997.9864563085441, 3149.75096277698),
    @(45515.99999999999, 2058, 992.4107314518574, 3110.077998846202),
    @(45522.99999999999, 2081, 1036.114147150108, 3146.153845713038),
    @(45529.99999999999, 2105, 1010.231981558529, 3198.458064691179),
    @(45550.99999999999, 2175, 1166.522518698401, 3286.166734689829),
    @(45557.99999999999, 2198, 1105.624009539906, 3374.8731351469),
    @(45564.99999999999, 2222, 1187.554697442598, 3394.72196755311),
    @(45571.99999999999, 2245, 1243.894812128134, 3318.392272979681),
    @(45578.99999999999, 2269, 1152.337558714035, 3334.153796713797),
    @(45585.99999999999, 2292, 1259.572196455808, 3388.424546240827),
    @(45599.99999999999, 2339, 1228.674096850313, 3425.016141650075),
    @(45606.99999999999, 2362, 1203.963991641865, 3430.309079789887),
    @(45613.99999999999, 2386, 1307.201296102289, 3460.535105368263),
    @(45620.99999999999, 2409, 1301.558732077362, 3450.645195299838),
    @(45627.99999999999, 2433, 1350.201404348636, 3483.612581156778),
    @(45634.99999999999, 2456, 1387.402094021114, 3445.224576080674),
    @(45641.99999999999, 2479, 1374.292190866729, 3515.788632116222),
    @(45648.99999999999, 2503, 1454.190600742586, 3637.551780780529),
    @(45655.99999999999, 2526, 1450.00690603299, 3566.89641152331)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $r = $data[$i]
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
}

# --- Match formatting of existing sheets: bold/bordered/centered header style,
#     and date-number-format for the "ds" column ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A86").PasteSpecial(-4122)

$excel.CutCopyMode = 0
$wsForecast.Range("A1").Select()

Write-Host "Done"
